$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new predidx (D), new pred_name (E)
$updates = @(
    @{Row=3;   D="[1, 0, 0, 1, 0, 0, 0]"; E="['Normal', 'ParamViolation']"},
    @{Row=12;  D="[1, 0, 0, 0, 0, 0, 0]"; E="['Normal']"},
    @{Row=15;  D="[0, 0, 0, 1, 0, 0, 0]"; E="['ParamViolation']"},
    @{Row=16;  D="[1, 0, 0, 0, 1, 0, 0]"; E="['Normal', 'RegulationViolation']"},
    @{Row=24;  D="[0, 0, 0, 0, 0, 0, 0]"; E="[]"},
    @{Row=26;  D="[0, 0, 0, 0, 0, 0, 0]"; E="[]"},
    @{Row=29;  D="[0, 0, 0, 1, 0, 0, 1]"; E="['ParamViolation', 'SoftwareFault']"},
    @{Row=31;  D="[1, 0, 0, 0, 0, 0, 1]"; E="['Normal', 'SoftwareFault']"},
    @{Row=36;  D="[1, 1, 0, 0, 0, 0, 0]"; E="['Normal', 'SurroundingEnvironment']"},
    @{Row=46;  D="[1, 0, 1, 0, 0, 0, 0]"; E="['Normal', 'HardwareFault']"},
    @{Row=54;  D="[0, 0, 0, 0, 0, 0, 0]"; E="[]"},
    @{Row=71;  D="[1, 0, 0, 0, 0, 0, 0]"; E="['Normal']"},
    @{Row=74;  D="[1, 0, 0, 0, 0, 0, 1]"; E="['Normal', 'SoftwareFault']"},
    @{Row=80;  D="[1, 0, 1, 0, 0, 0, 0]"; E="['Normal', 'HardwareFault']"},
    @{Row=81;  D="[1, 0, 1, 0, 0, 0, 0]"; E="['Normal', 'HardwareFault']"},
    @{Row=92;  D="[1, 0, 1, 0, 0, 0, 1]"; E="['Normal', 'HardwareFault', 'SoftwareFault']"},
    @{Row=93;  D="[1, 0, 1, 0, 0, 0, 1]"; E="['Normal', 'HardwareFault', 'SoftwareFault']"},
    @{Row=109; D="[1, 0, 0, 0, 0, 0, 0]"; E="['Normal']"},
    @{Row=113; D="[1, 0, 0, 0, 0, 0, 0]"; E="['Normal']"},
    @{Row=116; D="[1, 0, 0, 0, 0, 0, 0]"; E="['Normal']"}
)

foreach ($u in $updates) {
    $ws.Range("D$($u.Row)").Value = $u.D
    $ws.Range("E$($u.Row)").Value = $u.E
}
